$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell contents to reflect the new activity descriptions.
# Order matters: new unique strings are appended to the shared string table
# in the order they are first assigned, so we set cells in the sequence
# that reproduces the target shared string ordering (existing strings
# first, then the three brand-new strings in the order: B4, B5, B3).
$ws.Range("B6").Value = "Construir  Interfaces dos servicos Vacinacao, Exames e Cirurgia"
$ws.Range("B7").Value = "Construir Classe Controller   Vacinacao, Exames e Cirurgia"
$ws.Range("B8").Value = "Construitr classes DAO e fazer Conexao BD das classes  Vacinacao, Exames e Cirurgia"
$ws.Range("B4").Value = "Construir Classe Controller  Animal,Vacina , Cliente  e veterinaria"
$ws.Range("B5").Value = "Construitr classes DAO e fazer Conexao BD das classes   Animal,Vacina , Cliente e veterinaria"
$ws.Range("B3").Value = "Construir  Interfaces de cadastro e consulta do Animal,Vacina , Cliente e veterinaria  (interface update)"

# Adjust column widths (values chosen so the engine's internal pixel
# rounding lands on the stored width closest to the target: 92.5 / 6.667)
$ws.Columns.Item(2).ColumnWidth = 91.66666666666667
$ws.Columns.Item(3).ColumnWidth = 5.833333333333333
